$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-17 Saturday", 2) | Out-Null
$d.Content.Find.Execute("20×22=", $true, $false, $false, $false, $false, $true, 1, $false, "59×13=", 2) | Out-Null
$d.Content.Find.Execute("49×54=", $true, $false, $false, $false, $false, $true, 1, $false, "60×85=", 2) | Out-Null
$d.Content.Find.Execute("86×93=", $true, $false, $false, $false, $false, $true, 1, $false, "99×62=", 2) | Out-Null
$d.Content.Find.Execute("42×15=", $true, $false, $false, $false, $false, $true, 1, $false, "49×71=", 2) | Out-Null
$d.Content.Find.Execute("51×98=", $true, $false, $false, $false, $false, $true, 1, $false, "38×30=", 2) | Out-Null
$d.Content.Find.Execute("74×78=", $true, $false, $false, $false, $false, $true, 1, $false, "65×59=", 2) | Out-Null
$d.Content.Find.Execute("59×42=", $true, $false, $false, $false, $false, $true, 1, $false, "63×25=", 2) | Out-Null
$d.Content.Find.Execute("86×33=", $true, $false, $false, $false, $false, $true, 1, $false, "39×28=", 2) | Out-Null
$d.Content.Find.Execute("13×35=", $true, $false, $false, $false, $false, $true, 1, $false, "10×77=", 2) | Out-Null
$d.Content.Find.Execute("29×56=", $true, $false, $false, $false, $false, $true, 1, $false, "27×66=", 2) | Out-Null
$d.Content.Find.Execute("85×68=", $true, $false, $false, $false, $false, $true, 1, $false, "36×29=", 2) | Out-Null
$d.Content.Find.Execute("57×35=", $true, $false, $false, $false, $false, $true, 1, $false, "80×70=", 2) | Out-Null
$d.Content.Find.Execute("58×42=", $true, $false, $false, $false, $false, $true, 1, $false, "28×53=", 2) | Out-Null
$d.Content.Find.Execute("62×83=", $true, $false, $false, $false, $false, $true, 1, $false, "88×19=", 2) | Out-Null
$d.Content.Find.Execute("74×87=", $true, $false, $false, $false, $false, $true, 1, $false, "75×71=", 2) | Out-Null
$d.Content.Find.Execute("21×62=", $true, $false, $false, $false, $false, $true, 1, $false, "59×26=", 2) | Out-Null
$d.Content.Find.Execute("35×50=", $true, $false, $false, $false, $false, $true, 1, $false, "56×90=", 2) | Out-Null
$d.Content.Find.Execute("73×24=", $true, $false, $false, $false, $false, $true, 1, $false, "84×47=", 2) | Out-Null
$d.Content.Find.Execute("83×77=", $true, $false, $false, $false, $false, $true, 1, $false, "24×26=", 2) | Out-Null
$d.Content.Find.Execute("33×27=", $true, $false, $false, $false, $false, $true, 1, $false, "43×92=", 2) | Out-Null
$d.Content.Find.Execute("60×84=", $true, $false, $false, $false, $false, $true, 1, $false, "55×90=", 2) | Out-Null
$d.Content.Find.Execute("20×51=", $true, $false, $false, $false, $false, $true, 1, $false, "85×44=", 2) | Out-Null
$d.Content.Find.Execute("48×24=", $true, $false, $false, $false, $false, $true, 1, $false, "71×50=", 2) | Out-Null
$d.Content.Find.Execute("75×50=", $true, $false, $false, $false, $false, $true, 1, $false, "47×99=", 2) | Out-Null
$d.Content.Find.Execute("25×19=", $true, $false, $false, $false, $false, $true, 1, $false, "28×62=", 2) | Out-Null
$d.Content.Find.Execute("98×23=", $true, $false, $false, $false, $false, $true, 1, $false, "44×86=", 2) | Out-Null
$d.Content.Find.Execute("71×98=", $true, $false, $false, $false, $false, $true, 1, $false, "45×30=", 2) | Out-Null
$d.Content.Find.Execute("78×12=", $true, $false, $false, $false, $false, $true, 1, $false, "17×96=", 2) | Out-Null
$d.Content.Find.Execute("66×16=", $true, $false, $false, $false, $false, $true, 1, $false, "75×25=", 2) | Out-Null
$d.Content.Find.Execute("79×26=", $true, $false, $false, $false, $false, $true, 1, $false, "26×15=", 2) | Out-Null
$d.Content.Find.Execute("75×47=", $true, $false, $false, $false, $false, $true, 1, $false, "97×95=", 2) | Out-Null
$d.Content.Find.Execute("55×36=", $true, $false, $false, $false, $false, $true, 1, $false, "24×31=", 2) | Out-Null
$d.Content.Find.Execute("54×31=", $true, $false, $false, $false, $false, $true, 1, $false, "53×68=", 2) | Out-Null
$d.Content.Find.Execute("92×79=", $true, $false, $false, $false, $false, $true, 1, $false, "54×53=", 2) | Out-Null
$d.Content.Find.Execute("32×95=", $true, $false, $false, $false, $false, $true, 1, $false, "66×22=", 2) | Out-Null
$d.Content.Find.Execute("61×53=", $true, $false, $false, $false, $false, $true, 1, $false, "40×100=", 2) | Out-Null
$d.Content.Find.Execute("44×91=", $true, $false, $false, $false, $false, $true, 1, $false, "50×65=", 2) | Out-Null
$d.Content.Find.Execute("66×68=", $true, $false, $false, $false, $false, $true, 1, $false, "72×84=", 2) | Out-Null
$d.Content.Find.Execute("67×83=", $true, $false, $false, $false, $false, $true, 1, $false, "97×36=", 2) | Out-Null
$d.Content.Find.Execute("21×23=", $true, $false, $false, $false, $false, $true, 1, $false, "68×12=", 2) | Out-Null
$d.Content.Find.Execute("72×74=", $true, $false, $false, $false, $false, $true, 1, $false, "39×47=", 2) | Out-Null
$d.Content.Find.Execute("87×71=", $true, $false, $false, $false, $false, $true, 1, $false, "62×22=", 2) | Out-Null
$d.Content.Find.Execute("94×52=", $true, $false, $false, $false, $false, $true, 1, $false, "18×10=", 2) | Out-Null
$d.Content.Find.Execute("67×44=", $true, $false, $false, $false, $false, $true, 1, $false, "89×84=", 2) | Out-Null
$d.Content.Find.Execute("57×92=", $true, $false, $false, $false, $false, $true, 1, $false, "65×92=", 2) | Out-Null
$d.Content.Find.Execute("66×38=", $true, $false, $false, $false, $false, $true, 1, $false, "95×95=", 2) | Out-Null
$d.Content.Find.Execute("100×12=", $true, $false, $false, $false, $false, $true, 1, $false, "89×70=", 2) | Out-Null
$d.Content.Find.Execute("15×89=", $true, $false, $false, $false, $false, $true, 1, $false, "38×83=", 2) | Out-Null
$d.Content.Find.Execute("41×10=", $true, $false, $false, $false, $false, $true, 1, $false, "90×26=", 2) | Out-Null
$d.Content.Find.Execute("45×43=", $true, $false, $false, $false, $false, $true, 1, $false, "23×53=", 2) | Out-Null
$d.Content.Find.Execute("58×56=", $true, $false, $false, $false, $false, $true, 1, $false, "86×62=", 2) | Out-Null
$d.Content.Find.Execute("99×70=", $true, $false, $false, $false, $false, $true, 1, $false, "63×75=", 2) | Out-Null
$d.Content.Find.Execute("75×38=", $true, $false, $false, $false, $false, $true, 1, $false, "80×51=", 2) | Out-Null
$d.Content.Find.Execute("88×96=", $true, $false, $false, $false, $false, $true, 1, $false, "47×59=", 2) | Out-Null
$d.Content.Find.Execute("11×94=", $true, $false, $false, $false, $false, $true, 1, $false, "58×66=", 2) | Out-Null
$d.Content.Find.Execute("83×50=", $true, $false, $false, $false, $false, $true, 1, $false, "80×65=", 2) | Out-Null
$d.Content.Find.Execute("90×16=", $true, $false, $false, $false, $false, $true, 1, $false, "82×99=", 2) | Out-Null
$d.Content.Find.Execute("97×19=", $true, $false, $false, $false, $false, $true, 1, $false, "84×12=", 2) | Out-Null
$d.Content.Find.Execute("19×52=", $true, $false, $false, $false, $false, $true, 1, $false, "37×82=", 2) | Out-Null
$d.Content.Find.Execute("13×98=", $true, $false, $false, $false, $false, $true, 1, $false, "22×53=", 2) | Out-Null
$d.Content.Find.Execute("60×52=", $true, $false, $false, $false, $false, $true, 1, $false, "70×99=", 2) | Out-Null
$d.Content.Find.Execute("20×97=", $true, $false, $false, $false, $false, $true, 1, $false, "53×10=", 2) | Out-Null
$d.Content.Find.Execute("40×28=", $true, $false, $false, $false, $false, $true, 1, $false, "44×78=", 2) | Out-Null
$d.Content.Find.Execute("73×71=", $true, $false, $false, $false, $false, $true, 1, $false, "84×50=", 2) | Out-Null
$d.Content.Find.Execute("46×22=", $true, $false, $false, $false, $false, $true, 1, $false, "49×98=", 2) | Out-Null
$d.Content.Find.Execute("98×93=", $true, $false, $false, $false, $false, $true, 1, $false, "81×15=", 2) | Out-Null
$d.Content.Find.Execute("58×22=", $true, $false, $false, $false, $false, $true, 1, $false, "97×62=", 2) | Out-Null
$d.Content.Find.Execute("58×93=", $true, $false, $false, $false, $false, $true, 1, $false, "23×63=", 2) | Out-Null
$d.Content.Find.Execute("16×84=", $true, $false, $false, $false, $false, $true, 1, $false, "32×99=", 2) | Out-Null
$d.Content.Find.Execute("73×90=", $true, $false, $false, $false, $false, $true, 1, $false, "85×76=", 2) | Out-Null
$d.Content.Find.Execute("55×30=", $true, $false, $false, $false, $false, $true, 1, $false, "78×50=", 2) | Out-Null
$d.Content.Find.Execute("42×63=", $true, $false, $false, $false, $false, $true, 1, $false, "93×46=", 2) | Out-Null
$d.Content.Find.Execute("11×29=", $true, $false, $false, $false, $false, $true, 1, $false, "26×90=", 2) | Out-Null
$d.Content.Find.Execute("69×30=", $true, $false, $false, $false, $false, $true, 1, $false, "25×100=", 2) | Out-Null
$d.Content.Find.Execute("81×39=", $true, $false, $false, $false, $false, $true, 1, $false, "54×35=", 2) | Out-Null
$d.Content.Find.Execute("63×50=", $true, $false, $false, $false, $false, $true, 1, $false, "77×90=", 2) | Out-Null
$d.Content.Find.Execute("18×100=", $true, $false, $false, $false, $false, $true, 1, $false, "43×77=", 2) | Out-Null
$d.Content.Find.Execute("100×32=", $true, $false, $false, $false, $false, $true, 1, $false, "74×97=", 2) | Out-Null
$d.Content.Find.Execute("85×15=", $true, $false, $false, $false, $false, $true, 1, $false, "24×67=", 2) | Out-Null
$d.Content.Find.Execute("85×41=", $true, $false, $false, $false, $false, $true, 1, $false, "92×57=", 2) | Out-Null
$d.Content.Find.Execute("32×20=", $true, $false, $false, $false, $false, $true, 1, $false, "37×77=", 2) | Out-Null
$d.Content.Find.Execute("97×63=", $true, $false, $false, $false, $false, $true, 1, $false, "63×95=", 2) | Out-Null
$d.Content.Find.Execute("55×62=", $true, $false, $false, $false, $false, $true, 1, $false, "67×25=", 2) | Out-Null
$d.Content.Find.Execute("51×12=", $true, $false, $false, $false, $false, $true, 1, $false, "10×14=", 2) | Out-Null
$d.Content.Find.Execute("11×64=", $true, $false, $false, $false, $false, $true, 1, $false, "69×63=", 2) | Out-Null
$d.Content.Find.Execute("95×55=", $true, $false, $false, $false, $false, $true, 1, $false, "98×22=", 2) | Out-Null
$d.Content.Find.Execute("58×83=", $true, $false, $false, $false, $false, $true, 1, $false, "91×78=", 2) | Out-Null
$d.Content.Find.Execute("51×20=", $true, $false, $false, $false, $false, $true, 1, $false, "100×86=", 2) | Out-Null
$d.Content.Find.Execute("92×32=", $true, $false, $false, $false, $false, $true, 1, $false, "15×48=", 2) | Out-Null
$d.Content.Find.Execute("25×21=", $true, $false, $false, $false, $false, $true, 1, $false, "96×41=", 2) | Out-Null
$d.Content.Find.Execute("82×29=", $true, $false, $false, $false, $false, $true, 1, $false, "73×66=", 2) | Out-Null
$d.Content.Find.Execute("39×39=", $true, $false, $false, $false, $false, $true, 1, $false, "100×36=", 2) | Out-Null
$d.Content.Find.Execute("48×53=", $true, $false, $false, $false, $false, $true, 1, $false, "45×74=", 2) | Out-Null
$d.Content.Find.Execute("96×75=", $true, $false, $false, $false, $false, $true, 1, $false, "89×56=", 2) | Out-Null
$d.Content.Find.Execute("70×88=", $true, $false, $false, $false, $false, $true, 1, $false, "81×45=", 2) | Out-Null
$d.Content.Find.Execute("37×97=", $true, $false, $false, $false, $false, $true, 1, $false, "76×99=", 2) | Out-Null
$d.Content.Find.Execute("18×53=", $true, $false, $false, $false, $false, $true, 1, $false, "100×70=", 2) | Out-Null
$d.Content.Find.Execute("11×16=", $true, $false, $false, $false, $false, $true, 1, $false, "32×77=", 2) | Out-Null
$d.Content.Find.Execute("61×41=", $true, $false, $false, $false, $false, $true, 1, $false, "17×10=", 2) | Out-Null
$d.Content.Find.Execute("32×80=", $true, $false, $false, $false, $false, $true, 1, $false, "26×64=", 2) | Out-Null
